# Tilfoejede Lasse til to controller klassers
#
# Adds the text "Lasse" to the (currently empty) bullet paragraph that
# follows the "Logik (StaffController)" heading, and to the bullet
# paragraph that follows the "Logik (MainController)" heading (which
# previously only contained a manual line break).
#
# The new runs must carry the same run formatting (rFonts theme fonts +
# da-DK language) already used by every other name in the list, so
# rather than typing text in directly (which would pick up default
# formatting) we copy the FormattedText from an existing "Lasse" entry
# (the one under "Projektleder") and paste it into the target spots.

$d = $word.ActiveDocument

function Find-HeadingParagraph($doc, $headingText) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $headingText) {
            return $p
        }
    }
    return $null
}

# Donor text: the existing "Lasse" run under the "Projektleder" heading,
# which already has the correct rFonts(asciiTheme/hAnsiTheme=majorHAnsi)
# + lang=da-DK run formatting that every roster entry uses.
$projektlederHeading = Find-HeadingParagraph $d "Projektleder"
$lasseDonor = $projektlederHeading.Next()
$donorRange = $d.Range($lasseDonor.Range.Start, $lasseDonor.Range.End - 1)
$lasseFormatted = $donorRange.FormattedText

# --- 1) "Logik (StaffController)" -> its (empty) bullet paragraph gets "Lasse" ---
$staffHeading = Find-HeadingParagraph $d "Logik (StaffController)"
$staffListPara = $staffHeading.Next()
$staffInsertPos = $staffListPara.Range.Start
$d.Range($staffInsertPos, $staffInsertPos).FormattedText = $lasseFormatted

# --- 2) "Logik (MainController)" -> its bullet paragraph (which only had a
#        manual line break + the _GoBack bookmark) gets "Lasse" inserted
#        before the break, with the _GoBack bookmark ending up between the
#        new text and the break (matching the target layout). ---
$mainHeading = Find-HeadingParagraph $d "Logik (MainController)"
$mainListPara = $mainHeading.Next()
$mainInsertPos = $mainListPara.Range.Start
$d.Range($mainInsertPos, $mainInsertPos).FormattedText = $lasseFormatted

# Re-seat the _GoBack bookmark so it sits right after the freshly inserted
# "Lasse" text (i.e. before the line-break run), rather than at the very
# end of the paragraph where it landed after the insert shifted it along.
$goBack = $d.Bookmarks.Item("_GoBack")
$newBookmarkPos = $mainInsertPos + 5
$d.Bookmarks.Add("_GoBack", $d.Range($newBookmarkPos, $newBookmarkPos))

Write-Output "Added 'Lasse' under StaffController and MainController."
